$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing note text on row 17 (D17) ---
$ws.Range("D17").Value = "Adding in BA student to widen the model training pool (reversed..)"

# --- Update hours worked on row 17 (B17): 1 -> 2 ---
$ws.Range("B17").Value = 2

# --- Add new row 18: new date entry with hours, notes and a link ---
$ws.Range("A18").Value = 45431
$ws.Range("A18").NumberFormat = $ws.Range("A17").NumberFormat
$ws.Range("B18").Value = 4

# F18 gets the rich-text note first (so it becomes shared-string index 20),
# then D18 gets the plain note (shared-string index 21) - matches authoring order.
$ws.Range("F18").Value = "Look into best grade student receives"
$chars = $ws.Range("F18").Characters(11, 4)
$chars.Font.Bold = $true
$rest = $ws.Range("F18").Characters(15, 24)
$rest.Font.Bold = $false
$rest.Font.Name = "Aptos Narrow"
$rest.Font.Size = 11

$ws.Range("D18").Value = "Creating a new dataframe to save only newest grades student recieves, weights are based on rows in GBM"

# Row 18 wraps to a taller row, matching the multi-line note content.
$ws.Rows.Item(18).RowHeight = 45
